# tdf#141309 tdf#142215 - OOXML import: fix double conversion in autofilter
# resulting in a missing selection of filtered time values and numbers
# ending with zeroes.
#
# This test document gets a new data row (time value 0.5 / "c") appended,
# the autofilter range is extended to cover it, and the filter criteria on
# column A gain the matching "0.500" value so the new row stays visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 8): time value 0.5 in column A, "c" in column B.
$ws.Range("A8").Value = 0.5
$ws.Range("B8").Value = "c"

# Drop the existing autofilter so it can be reapplied over the widened
# range (A1:B7 -> A1:B8) together with the extra "0.500" filter criterion.
$ws.AutoFilterMode = $false
$ws.Range("A1:B8").AutoFilter(1, @("0.046", "0.500", "0.516"), 7)

# Keep the filter database defined name in sync with the new range.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Munka1!`$A`$1:`$B`$8"

# Match the author's final selection.
$ws.Range("C7").Select()
